$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: "Must be in format "wXXXXXXX"" -> split "wXXXXXXX" into its own
# run, wrapped in proofErr spellStart/spellEnd markers (as Word does for
# words it flags during spell-check), splitting the original single run
# into three runs.
# ---------------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*wXXXXXXX*") {
        $pStart = $p.Range.Start
        $pEnd = $p.Range.End
        # Range covering just the paragraph's text, excluding the trailing
        # paragraph mark.
        $textRange = $d.Range($pStart, $pEnd - 1)
        $textRange.Text = ""

        $collapsed = $d.Range($pStart, $pStart)
        $frag = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>Must be in format &#8220;</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>wXXXXXXX</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>&#8221;</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
        $collapsed.InsertXML($frag)
        break
    }
}

# ---------------------------------------------------------------------------
# Change 2: "Max 4 characters (XX.XX)" -> "Max 3 characters (XX.X)", keeping
# the paragraph's existing four-run split (only the text inside each run
# changes, plus the now-required/no-longer-required xml:space="preserve").
# ---------------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Max 4 characters*") {
        $pStart = $p.Range.Start
        $pEnd = $p.Range.End
        $textRange = $d.Range($pStart, $pEnd - 1)
        $textRange.Text = ""

        $collapsed = $d.Range($pStart, $pStart)
        $frag = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">Max </w:t></w:r><w:r w:rsidR="00FD2D23"><w:t xml:space="preserve">3 </w:t></w:r><w:r><w:t>characters</w:t></w:r><w:r w:rsidR="00DA72E8"><w:t xml:space="preserve"> (XX.X)</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
        $collapsed.InsertXML($frag)
        break
    }
}
